# Apply the cryptos list update (Wed Dec 27 16:41:05 UTC 2023 GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.233.62'
$ws.Range("E2").Value = '  +1.61%  '
$ws.Range("D3").Value = '2.370.89'
$ws.Range("E3").Value = '  +6.62%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = "'311.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.47%  '
$ws.Range("D6").Value = "'107.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.84%  '
$ws.Range("D7").Value = "'0.641"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.11%  '
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").Value = "'0.633"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.83%  '
$ws.Range("D10").Value = "'43.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.13%  '
$ws.Range("D11").Value = "'0.0941"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("D12").Value = "'9.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = "'1.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.89%  '
$ws.Range("D14").Value = "'16.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.24%  '
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '2.728.30'
$ws.Range("E16").Value = '  +6.54%  '
$ws.Range("D17").Value = '2.357.22'
$ws.Range("E17").Value = '  +5.18%  '
$ws.Range("D18").Value = '43.244.08'
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").Value = "'0.0000109"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").Value = "'7.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.46%  '
$ws.Range("D21").Value = "'75.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.50%  '
$ws.Range("D22").Value = "'3.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("D23").Value = "'2.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.01%  '
$ws.Range("D24").Value = "'253.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'8.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.49%  '
$ws.Range("D26").Value = "'12.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.35%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = "'39.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = "'2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.52%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = "'22.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.70%  '
$ws.Range("D31").Value = "'173.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("D33").Value = "'0.0911"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.06%  '
$ws.Range("D34").Value = "'5.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("D35").Value = "'5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("D36").Value = "'0.132"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.61%  '
$ws.Range("D37").Value = "'0.0379"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").Value = "'4.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.41%  '
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("D40").Value = "'2.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.51%  '
$ws.Range("D41").Value = "'1.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +16.24%  '
$ws.Range("D42").Value = "'72.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").Value = "'0.234"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").Value = "'12.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.69%  '
$ws.Range("D46").Value = "'5.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.32%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'9.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.52%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = "'112.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.08%  '
$ws.Range("E49").Value = '  -1.78%  '
$ws.Range("D50").Value = "'0.0998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("D51").Value = '1.498.08'
$ws.Range("E51").Value = '  +4.33%  '
